# add genders in batumi
# Add a new column S (year 2023) to the trade data sheet, mirroring
# the existing column R (year 2022) formatting, and update the
# selection/view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column S data ----------------------------------------------------
# Row 1: header continuation (merged title cell area), keep same style as R1
$ws.Range("S1").Value = $null

# Row 3: year header
$ws.Range("S3").Value = 2023

# Data rows 4-14
$ws.Range("S4").Value = 200
$ws.Range("S5").Value = 45.8
$ws.Range("S6").Value = 2005
$ws.Range("S7").Value = 1134
$ws.Range("S8").Value = 526.9
$ws.Range("S9").Value = 6.8
$ws.Range("S10").Value = 7.2
$ws.Range("S11").Value = 38
$ws.Range("S12").Value = 1.5
$ws.Range("S13").Value = 166
$ws.Range("S14").Value = 159.6

# --- Copy formatting from column R into column S for rows 1, 3-14 --------
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Extend the merged title cell A1:R1 to A1:S1 --------------------------
$ws.Range("A1:R1").UnMerge() | Out-Null
$ws.Range("A1:S1").Merge() | Out-Null

# --- Update the view: scroll so column I is the top-left, select S3:S14 --
$ws.Range("S3:S14").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1

# --- Match the resized application window recorded in the workbook view --
$excel.ActiveWindow.Width = 16395
$excel.ActiveWindow.Height = 9285
$excel.CutCopyMode = 0
